# Horarios actualizados Linea 141 - 196
# Updates schedule data across sheets LP1912, LP1912-215, 6203-6173
$wb = $excel.ActiveWorkbook

$newUpdateTime = "05:20:00"

# ---- Sheet 1: LP1912 ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Cells.Item(2, 1).Value = "Última actualización: " + $newUpdateTime
$ws1.Cells.Item(3, 1).Value = "Total filas: 31"

$rows1 = @(
    ,@(13, "05:20:00", "05:22", "14_ABASTO", 2, "LP1912")
    ,@(14, "04:01:01", "05:22", "23_HERNANDEZ", 81, "LP1912")
    ,@(15, "04:36:47", "05:34", "215B_EL PATO", 58, "LP1912")
    ,@(16, "04:01:01", "05:35", "215B_EL PATO", 94, "LP1912")
    ,@(17, "04:01:01", "05:41", "14_ABASTO", 100, "LP1912")
    ,@(18, "04:01:01", "05:46", "15_ABASTO", 105, "LP1912")
    ,@(19, "04:36:47", "06:04", "16_SANTA ANA", 88, "LP1912")
    ,@(20, "04:36:47", "06:11", "215A_EL PATO", 95, "LP1912")
    ,@(21, "05:20:00", "06:12", "215A_EL PATO", 52, "LP1912")
    ,@(22, "04:36:47", "06:14", "225_HARAS DEL SUR", 98, "LP1912")
    ,@(23, "04:36:47", "06:21", "26_HERNANDEZ", 105, "LP1912")
    ,@(24, "04:36:47", "06:27", "23_HERNANDEZ", 111, "LP1912")
    ,@(25, "04:36:47", "06:29", "86_EST CHICA-ESC AGRARIA", 113, "LP1912")
    ,@(26, "05:20:00", "06:30", "86_EST CHICA-ESC AGRARIA", 70, "LP1912")
    ,@(27, "04:36:47", "06:31", "16_SANTA ANA", 115, "LP1912")
    ,@(28, "04:51:28", "06:44", "225_C ROCA-H SUR", 113, "LP1912")
    ,@(29, "04:51:28", "06:46", "215C_EL PATO", 115, "LP1912")
    ,@(30, "05:20:00", "06:47", "215C_EL PATO", 87, "LP1912")
    ,@(31, "05:20:00", "07:00", "10_OLMOS", 100, "LP1912")
    ,@(32, "05:20:00", "07:00", "14_ABASTO", 100, "LP1912")
    ,@(33, "05:20:00", "07:05", "15_ABASTO", 105, "LP1912")
    ,@(34, "05:20:00", "07:07", "225_GOMEZ", 107, "LP1912")
    ,@(35, "05:20:00", "07:12", "215A_EL PATO", 112, "LP1912")
    ,@(36, "05:20:00", "07:16", "11_ETCHEVERRY", 116, "LP1912")
)

foreach ($row in $rows1) {
    $r = $row[0]
    $ws1.Cells.Item($r, 1).Value = $row[1]
    $ws1.Cells.Item($r, 2).Value = $row[2]
    $ws1.Cells.Item($r, 3).Value = $row[3]
    $ws1.Cells.Item($r, 4).Value = $row[4]
    $ws1.Cells.Item($r, 5).Value = $row[5]
}

# ---- Sheet 2: LP1912-215 ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2, 1).Value = "Última actualización: " + $newUpdateTime
$ws2.Cells.Item(3, 1).Value = "Total filas: 9"

$rows2 = @(
    ,@(11, "05:20:00", "06:12", "215A_EL PATO", 52, "LP1912")
    ,@(12, "04:51:28", "06:46", "215C_EL PATO", 115, "LP1912")
    ,@(13, "05:20:00", "06:47", "215C_EL PATO", 87, "LP1912")
    ,@(14, "05:20:00", "07:12", "215A_EL PATO", 112, "LP1912")
)

foreach ($row in $rows2) {
    $r = $row[0]
    $ws2.Cells.Item($r, 1).Value = $row[1]
    $ws2.Cells.Item($r, 2).Value = $row[2]
    $ws2.Cells.Item($r, 3).Value = $row[3]
    $ws2.Cells.Item($r, 4).Value = $row[4]
    $ws2.Cells.Item($r, 5).Value = $row[5]
}

# ---- Sheet 3: 6203-6173 ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2, 1).Value = "Última actualización: " + $newUpdateTime
$ws3.Cells.Item(3, 1).Value = "Total filas: 7"

$rows3 = @(
    ,@(12, "05:20:00", "07:00", "215B_LP-P MOR-1 Y 57", 100, "L6173")
)

foreach ($row in $rows3) {
    $r = $row[0]
    $ws3.Cells.Item($r, 1).Value = $row[1]
    $ws3.Cells.Item($r, 2).Value = $row[2]
    $ws3.Cells.Item($r, 3).Value = $row[3]
    $ws3.Cells.Item($r, 4).Value = $row[4]
    $ws3.Cells.Item($r, 5).Value = $row[5]
}
